# repull data, push all data, mean calculation
# Update the dSF (column F) values on Sheet1 to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -7
    7  = 0
    9  = -1
    10 = -2
    16 = 4
    17 = -2
    18 = 2
    21 = -4
    22 = -8
    29 = 4
    30 = -4
    31 = 0
    33 = -3
    36 = -10
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
